$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix the existing labels in column A (rows 2-11) with a letter code,
# e.g. "Alpha Quartz" -> "A-Alpha Quartz", "4000-50" -> "B-4000-50", etc.
$letters = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $letters.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = "$($letters[$i])-$($cell.Value2)"
}

$null = $ws.Range("A11").Select()
